# Command Strings List.xlsx edit
# - Added radius turning notes / fixed typo, updated waypoint precision example
# - Added a new "Forget waypoints" command row (erase all waypoints)
# - Removed some debug prints (typo fix "conenction" -> "connection")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix "conenction" -> "connection" typo in the GPS / temperature stream notes
$ws.Range("G5").Value = "Simply open the connection to receive data"
$ws.Range("G6").Value = "Simply open the connection to receive data"

# Update waypoint lat/long precision note and example from 4 to 6 decimal places
$ws.Range("G31").Value = " +lat is N, -lat is S; +long is E, -long is W               Values to 6 decimal places"
$ws.Range("H31").Value = "NW42.034534,-93.620369"

# Insert a new row above the existing "Disable" row (row 32) for the new
# "Forget waypoints" command, pushing Disable down to row 33.
$ws.Rows("32").Insert()

# The newly inserted blank row inherits formatting from the row above (row 31);
# clear the formatting on F:H so they match the plain (unstyled) target cells.
$ws.Range("F32:H32").ClearFormats()

$ws.Range("B32").Value = "Forget waypoints"
$ws.Range("D32").Value = "N"
$ws.Range("E32").Value = "F"
$ws.Range("F32").Value = "N/A"
$ws.Range("G32").Value = "Erases (forgets) all waypoints"
$ws.Range("H32").Value = "NF"

# Restore the view's active selection (it had moved to H32 before the edit).
$ws.Range("G6").Select()
